# Update workbook for the new IG build (release-notes.md f80ed2bb9e1dd81abc71d13817b8a44a756cee80):
#  - bump the published Version / Status / Date / Contact metadata
#  - swap the two "Mapping" columns (AK/AL) on the Elements sheet so that
#    "Spécification métier vers l'extension AsLieuDit" comes before "RIM Mapping"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metadata sheet: update the published Version / Status / Date / Contact
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Cells.Item(3, 2).Value2  = "0.4.0-snapshot-1"                 # Version
$meta.Cells.Item(6, 2).Value2  = "draft"                            # Status
$meta.Cells.Item(8, 2).Value2  = "2024-05-23T12:16:26+00:00"        # Date
$meta.Cells.Item(10, 2).Value2 = "ANS (https://esante.gouv.fr)"     # Contact

# ---------------------------------------------------------------------------
# Elements sheet: swap columns AK (37) and AL (38) - header + all data rows -
# so the "Spécification métier" mapping column now precedes "RIM Mapping".
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$lastRow = $elements.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elements.Cells.Item($r, 37)
    $alCell = $elements.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value2 = $alVal
    $alCell.Value2 = $akVal
}

# Column widths follow the content: AK becomes the wide "Spécification
# métier" column, AL becomes the narrower "RIM Mapping" column.
$elements.Columns.Item(37).ColumnWidth = 59.90234375
$elements.Columns.Item(38).ColumnWidth = 24.98046875
